# Apply the "Updated cryptos list" data refresh (prices / 1h volume deltas),
# including the OKB / Stacks row swap at rows 44-45.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "69.177.26"
$ws.Cells.Item(2, 5).Formula = "  +1.36%  "
$ws.Cells.Item(3, 4).Formula = "3.748.36"
$ws.Cells.Item(3, 5).Formula = "  +0.70%  "
$ws.Cells.Item(4, 5).Formula = "  -0.06%  "
$ws.Cells.Item(5, 4).Formula = "'602.94"
$ws.Cells.Item(5, 5).Formula = "  +0.70%  "
$ws.Cells.Item(6, 4).Formula = "'168.48"
$ws.Cells.Item(6, 5).Formula = "  +0.34%  "
$ws.Cells.Item(7, 4).Formula = "3.746.92"
$ws.Cells.Item(7, 5).Formula = "  +0.71%  "
$ws.Cells.Item(8, 5).Formula = "  -0.06%  "
$ws.Cells.Item(9, 4).Formula = "'0.542"
$ws.Cells.Item(9, 5).Formula = "  +1.30%  "
$ws.Cells.Item(10, 5).Formula = "  -0.19%  "
$ws.Cells.Item(11, 4).Formula = "'6.45"
$ws.Cells.Item(11, 5).Formula = "  +3.32%  "
$ws.Cells.Item(12, 5).Formula = "  -0.19%  "
$ws.Cells.Item(13, 4).Formula = "'38.22"
$ws.Cells.Item(13, 5).Formula = "  -0.48%  "
$ws.Cells.Item(14, 5).Formula = "  +0.86%  "
$ws.Cells.Item(15, 4).Formula = "4.381.62"
$ws.Cells.Item(15, 5).Formula = "  +0.71%  "
$ws.Cells.Item(16, 4).Formula = "3.768.94"
$ws.Cells.Item(16, 5).Formula = "  +1.29%  "
$ws.Cells.Item(17, 4).Formula = "69.185.12"
$ws.Cells.Item(17, 5).Formula = "  +1.36%  "
$ws.Cells.Item(18, 4).Formula = "'7.33"
$ws.Cells.Item(18, 5).Formula = "  -0.10%  "
$ws.Cells.Item(19, 4).Formula = "'0.114"
$ws.Cells.Item(19, 5).Formula = "  -1.28%  "
$ws.Cells.Item(20, 4).Formula = "'17.14"
$ws.Cells.Item(20, 5).Formula = "  -0.56%  "
$ws.Cells.Item(21, 4).Formula = "'11.01"
$ws.Cells.Item(21, 5).Formula = "  +18.80%  "
$ws.Cells.Item(22, 4).Formula = "'493.91"
$ws.Cells.Item(22, 5).Formula = "  +0.62%  "
$ws.Cells.Item(23, 4).Formula = "'0.724"
$ws.Cells.Item(23, 5).Formula = "  -0.12%  "
$ws.Cells.Item(24, 5).Formula = "  +5.93%  "
$ws.Cells.Item(25, 4).Formula = "'84.93"
$ws.Cells.Item(25, 5).Formula = "  -0.03%  "
$ws.Cells.Item(26, 4).Formula = "'2.31"
$ws.Cells.Item(26, 5).Formula = "  -0.32%  "
$ws.Cells.Item(27, 4).Formula = "'12.37"
$ws.Cells.Item(27, 5).Formula = "  +0.34%  "
$ws.Cells.Item(28, 4).Formula = "'10.16"
$ws.Cells.Item(28, 5).Formula = "  +0.38%  "
$ws.Cells.Item(29, 5).Formula = "  -0.13%  "
$ws.Cells.Item(30, 4).Formula = "'3.00"
$ws.Cells.Item(30, 5).Formula = "  +2.49%  "
$ws.Cells.Item(31, 4).Formula = "'2.49"
$ws.Cells.Item(31, 5).Formula = "  +4.62%  "
$ws.Cells.Item(32, 4).Formula = "'8.00"
$ws.Cells.Item(32, 5).Formula = "  +1.46%  "
$ws.Cells.Item(33, 4).Formula = "'31.58"
$ws.Cells.Item(33, 5).Formula = "  -0.09%  "
$ws.Cells.Item(34, 4).Formula = "3.896.39"
$ws.Cells.Item(34, 5).Formula = "  +0.73%  "
$ws.Cells.Item(35, 5).Formula = "  +0.27%  "
$ws.Cells.Item(36, 4).Formula = "3.686.42"
$ws.Cells.Item(36, 5).Formula = "  +0.51%  "
$ws.Cells.Item(37, 4).Formula = "'1.00"
$ws.Cells.Item(37, 5).Formula = "  -0.11%  "
$ws.Cells.Item(38, 5).Formula = "  +1.32%  "
$ws.Cells.Item(39, 4).Formula = "'5.87"
$ws.Cells.Item(39, 5).Formula = "  +0.45%  "
$ws.Cells.Item(40, 5).Formula = "  +0.71%  "
$ws.Cells.Item(41, 4).Formula = "'0.324"
$ws.Cells.Item(41, 5).Formula = "  -0.08%  "
$ws.Cells.Item(42, 4).Formula = "'2.97"
$ws.Cells.Item(42, 5).Formula = "  +4.01%  "
$ws.Cells.Item(43, 4).Formula = "'432.69"
$ws.Cells.Item(43, 5).Formula = "  +0.19%  "
$ws.Cells.Item(44, 2).Formula = "Stacks"
$ws.Cells.Item(44, 3).Formula = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(44, 4).Formula = "'1.99"
$ws.Cells.Item(44, 5).Formula = "  +1.31%  "
$ws.Cells.Item(45, 2).Formula = "OKB"
$ws.Cells.Item(45, 3).Formula = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(45, 4).Formula = "'48.52"
$ws.Cells.Item(45, 5).Formula = "  -0.45%  "
$ws.Cells.Item(46, 5).Formula = "  +0.75%  "
$ws.Cells.Item(47, 5).Formula = "  -0.04%  "
$ws.Cells.Item(48, 4).Formula = "'40.54"
$ws.Cells.Item(48, 5).Formula = "  +0.14%  "
$ws.Cells.Item(49, 4).Formula = "'141.32"
$ws.Cells.Item(49, 5).Formula = "  -0.04%  "
$ws.Cells.Item(50, 4).Formula = "2.792.88"
$ws.Cells.Item(50, 5).Formula = "  +1.18%  "
$ws.Cells.Item(51, 5).Formula = "  +0.21%  "
